$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (the blank spacer row between the data rows and the totals rows).
# This shifts rows 7-11 up to rows 6-10 and updates borders/merges/formulas accordingly.
$ws.Rows.Item(6).Delete()

# Update the active selection to match the post-edit state recorded in the workbook.
$ws.Range("E12").Select()
